# Weekly Progress Report 4 Updated.
#
# This script applies four related changes to the "Progress made in
# Reporting Week" section and one formatting tweak to a figure:
#
#   1. The "Bansil Patel:" bullet becomes "Bansil Patel & Meet Patel:"
#      and gains a new sentence describing the combined work done.
#   2. The now-redundant standalone "Meet Patel:" bullet (in the same
#      sub-section) is removed entirely.
#   3. The previously-empty bullet under "GitHub Management for
#      Deployment and Test-Automation." gets new descriptive text.
#   4. The screenshot illustrating "Figure 4: Model Evaluation Summary"
#      is marked NoProofing (<w:noProof/>) like its sibling figures.

$d = $word.ActiveDocument

function Get-ParaText($para) {
    $t = $para.Range.Text
    if ($t.Length -gt 0 -and [int][char]$t[$t.Length - 1] -eq 13) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    return $t
}

function Find-ParaIndexExact($exactText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ((Get-ParaText $d.Paragraphs.Item($i)) -eq $exactText) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndexContains($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ((Get-ParaText $d.Paragraphs.Item($i)).Contains($substr)) {
            return $i
        }
    }
    return -1
}

function Insert-RunsXmlBefore($range, $bodyXml) {
    $wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $bodyXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($wrapped)
}

# ---------------------------------------------------------------------
# 3. Fill in the previously-empty bullet right after "GitHub Management
#    for Deployment and Test-Automation." with the new task summary.
#    (Do the lower-in-document edits first so paragraph indices found
#    via text search above it stay valid for the edits that follow.)
# ---------------------------------------------------------------------
$idxGithubMgmt = Find-ParaIndexContains("GitHub Management for Deployment and Test-Automation.")
$idxEmptyTask = $idxGithubMgmt + 1
$pEmptyTask = $d.Paragraphs.Item($idxEmptyTask)
$insEmptyTask = $d.Range($pEmptyTask.Range.Start, $pEmptyTask.Range.Start)
$taskXml = '<w:r><w:rPr><w:b w:val="0"/><w:bCs/></w:rPr><w:t xml:space="preserve">Applying initial </w:t></w:r>'
$taskXml = $taskXml + '<w:r><w:rPr><w:b w:val="0"/><w:bCs/></w:rPr><w:t>NN function and testing the results on the sample data.</w:t></w:r>'
Insert-RunsXmlBefore $insEmptyTask $taskXml

# ---------------------------------------------------------------------
# 2. Delete the standalone "Meet Patel: " bullet that immediately
#    follows the "... test automation." bullet.
# ---------------------------------------------------------------------
$idxTestAutomation = Find-ParaIndexContains("test automation.")
$idxMeetPatelOnly = $idxTestAutomation + 1
$d.Paragraphs.Item($idxMeetPatelOnly).Range.Delete()

# ---------------------------------------------------------------------
# 1. Turn the "Bansil Patel: " bullet (exact text, trailing space, no
#    further content) into "Bansil Patel & Meet Patel:" followed by a
#    description of the combined work.
# ---------------------------------------------------------------------
$idxBansil = Find-ParaIndexExact("Bansil Patel: ")
$pBansil = $d.Paragraphs.Item($idxBansil)
$startBansil = $pBansil.Range.Start
$endBansil = $pBansil.Range.End - 1   # exclude the paragraph mark
$contentRange = $d.Range($startBansil, $endBansil)
$contentRange.Delete()
$insBansil = $d.Range($startBansil, $startBansil)
$bansilXml = '<w:r><w:t>Bansil Patel</w:t></w:r>'
$bansilXml = $bansilXml + '<w:r><w:t xml:space="preserve"> &amp; Meet Patel</w:t></w:r>'
$bansilXml = $bansilXml + '<w:r><w:t>:</w:t></w:r>'
$bansilXml = $bansilXml + '<w:r><w:rPr><w:b w:val="0"/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$bansilXml = $bansilXml + '<w:r><w:rPr><w:b w:val="0"/><w:bCs/></w:rPr><w:t>Applying NN and Model Evaluation through Summary.</w:t></w:r>'
Insert-RunsXmlBefore $insBansil $bansilXml

# ---------------------------------------------------------------------
# 4. Mark the "Figure 4: Model Evaluation Summary" screenshot as
#    NoProofing, matching the other figures in the document. Identify
#    it by its rendered size (6.5in x 3.625in = 468pt x 261pt) rather
#    than a document-order index, since that is stable and unique.
# ---------------------------------------------------------------------
$targetShape = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $candidate = $d.InlineShapes.Item($i)
    if ([Math]::Round($candidate.Width, 1) -eq 468 -and [Math]::Round($candidate.Height, 1) -eq 261) {
        $targetShape = $candidate
        break
    }
}
if ($targetShape -eq $null) {
    # Fallback: the figure immediately preceding the "Figure 4" caption.
    $idxFigure4Caption = Find-ParaIndexContains("Figure 4: Model Evaluation Summary")
    $targetShape = $d.InlineShapes.Item($idxFigure4Caption - 1)
}
$targetShape.Range.NoProofing = 1

Write-Output "Edit complete."
